$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.07899133333333333
$ws.Range("H2").Value = 0.236974
$ws.Range("I2").Value = 0.08952820636815619
$ws.Range("J2").Value = 0.09988914039546951
$ws.Range("M2").Value = 7.004922666666666
$ws.Range("N2").Value = 21.014768
$ws.Range("O2").Value = 0.134029393318039
$ws.Range("P2").Value = 0.1464771065395205
$ws.Range("Q2").Value = 0.5533281813368889
$ws.Range("R2").Value = 4.979953632032
$ws.Range("S2").Value = 0.01199941118437617
$ws.Range("T2").Value = 0.01463147225984831
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.07899133333333333
$ws.Range("H3").Value = 0.236974
$ws.Range("I3").Value = 0.08952820636815619
$ws.Range("J3").Value = 0.09988914039546951
$ws.Range("O3").Value = 0.3796257919253833
$ws.Range("P3").Value = 0.4148827819958515
$ws.Range("Q3").Value = 1.567250614469333
$ws.Range("R3").Value = 14.105255530224
$ws.Range("S3").Value = 0.03398721624217043
$ws.Range("T3").Value = 0.04144228445844658
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.07899133333333333
$ws.Range("H4").Value = 0.236974
$ws.Range("I4").Value = 0.08952820636815619
$ws.Range("J4").Value = 0.09988914039546951
$ws.Range("M4").Value = 6.293636666666667
$ws.Range("N4").Value = 18.88091
$ws.Range("O4").Value = 0.12041993100245
$ws.Range("P4").Value = 0.1316036924905903
$ws.Range("Q4").Value = 0.4971427518155556
$ws.Range("R4").Value = 4.47428476634
$ws.Range("S4").Value = 0.01078098043362648
$ws.Range("T4").Value = 0.01314577971575477
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.07899133333333333
$ws.Range("H5").Value = 0.236974
$ws.Range("I5").Value = 0.08952820636815619
$ws.Range("J5").Value = 0.09988914039546951
$ws.Range("M5").Value = 13.324299
$ws.Range("N5").Value = 26.648598
$ws.Range("O5").Value = 0.2549418168249328
$ws.Range("P5").Value = 0.1857460205306503
$ws.Range("Q5").Value = 1.052504143742
$ws.Range("R5").Value = 6.315024862452
$ws.Range("S5").Value = 0.02282448358857526
$ws.Range("T5").Value = 0.01855401032268588
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.07899133333333333
$ws.Range("H6").Value = 0.236974
$ws.Range("I6").Value = 0.08952820636815619
$ws.Range("J6").Value = 0.09988914039546951
$ws.Range("M6").Value = 5.800427666666667
$ws.Range("N6").Value = 17.401283
$ws.Range("O6").Value = 0.110983066929195
$ws.Range("P6").Value = 0.1212903984433873
$ws.Range("Q6").Value = 0.4581835152935556
$ws.Range("R6").Value = 4.123651637641999
$ws.Range("S6").Value = 0.009936114919407858
$ws.Range("T6").Value = 0.01211559363873396
$ws.Range("G7").Value = 0.296947
$ws.Range("H7").Value = 0.890841
$ws.Range("I7").Value = 0.336557583908845
$ws.Range("J7").Value = 0.3755067717093034
$ws.Range("M7").Value = 7.004922666666666
$ws.Range("N7").Value = 21.014768
$ws.Range("O7").Value = 0.134029393318039
$ws.Range("P7").Value = 0.1464771065395205
$ws.Range("Q7").Value = 2.080090771098667
$ws.Range("R7").Value = 18.720816939888
$ws.Range("S7").Value = 0.04510860878788749
$ws.Range("T7").Value = 0.05500314540597502
$ws.Range("G8").Value = 0.296947
$ws.Range("H8").Value = 0.890841
$ws.Range("I8").Value = 0.336557583908845
$ws.Range("J8").Value = 0.3755067717093034
$ws.Range("O8").Value = 0.3796257919253833
$ws.Range("P8").Value = 0.4148827819958515
$ws.Range("Q8").Value = 5.891663662024
$ws.Range("R8").Value = 53.024972958216
$ws.Range("S8").Value = 0.1277659393198889
$ws.Range("T8").Value = 0.1557912941050369
$ws.Range("G9").Value = 0.296947
$ws.Range("H9").Value = 0.890841
$ws.Range("I9").Value = 0.336557583908845
$ws.Range("J9").Value = 0.3755067717093034
$ws.Range("M9").Value = 6.293636666666667
$ws.Range("N9").Value = 18.88091
$ws.Range("O9").Value = 0.12041993100245
$ws.Range("P9").Value = 0.1316036924905903
$ws.Range("Q9").Value = 1.868876527256667
$ws.Range("R9").Value = 16.81988874531
$ws.Range("S9").Value = 0.04052824103265441
$ws.Range("T9").Value = 0.04941807771216546
$ws.Range("G10").Value = 0.296947
$ws.Range("H10").Value = 0.890841
$ws.Range("I10").Value = 0.336557583908845
$ws.Range("J10").Value = 0.3755067717093034
$ws.Range("M10").Value = 13.324299
$ws.Range("N10").Value = 26.648598
$ws.Range("O10").Value = 0.2549418168249328
$ws.Range("P10").Value = 0.1857460205306503
$ws.Range("Q10").Value = 3.956610615153
$ws.Range("R10").Value = 23.739663690918
$ws.Range("S10").Value = 0.08580260190793072
$ws.Range("T10").Value = 0.06974888852731446
$ws.Range("G11").Value = 0.296947
$ws.Range("H11").Value = 0.890841
$ws.Range("I11").Value = 0.336557583908845
$ws.Range("J11").Value = 0.3755067717093034
$ws.Range("M11").Value = 5.800427666666667
$ws.Range("N11").Value = 17.401283
$ws.Range("O11").Value = 0.110983066929195
$ws.Range("P11").Value = 0.1212903984433873
$ws.Range("Q11").Value = 1.722419594333667
$ws.Range("R11").Value = 15.501776349003
$ws.Range("S11").Value = 0.0373521928604835
$ws.Range("T11").Value = 0.04554536595881149
$ws.Range("G12").Value = 0.27455
$ws.Range("H12").Value = 0.5491
$ws.Range("I12").Value = 0.3111729859610415
$ws.Range("J12").Value = 0.2314563074056745
$ws.Range("M12").Value = 7.004922666666666
$ws.Range("N12").Value = 21.014768
$ws.Range("O12").Value = 0.134029393318039
$ws.Range("P12").Value = 0.1464771065395205
$ws.Range("Q12").Value = 1.923201518133333
$ws.Range("R12").Value = 11.5392091088
$ws.Range("S12").Value = 0.04170632652532105
$ws.Range("T12").Value = 0.03390305019910499
$ws.Range("G13").Value = 0.27455
$ws.Range("H13").Value = 0.5491
$ws.Range("I13").Value = 0.3111729859610415
$ws.Range("J13").Value = 0.2314563074056745
$ws.Range("O13").Value = 0.3796257919253833
$ws.Range("P13").Value = 0.4148827819958515
$ws.Range("Q13").Value = 5.447289443600001
$ws.Range("R13").Value = 32.6837366616
$ws.Range("S13").Value = 0.1181292912212466
$ws.Range("T13").Value = 0.09602723672695325
$ws.Range("G14").Value = 0.27455
$ws.Range("H14").Value = 0.5491
$ws.Range("I14").Value = 0.3111729859610415
$ws.Range("J14").Value = 0.2314563074056745
$ws.Range("M14").Value = 6.293636666666667
$ws.Range("N14").Value = 18.88091
$ws.Range("O14").Value = 0.12041993100245
$ws.Range("P14").Value = 0.1316036924905903
$ws.Range("Q14").Value = 1.727917946833333
$ws.Range("R14").Value = 10.367507681
$ws.Range("S14").Value = 0.03747142949925498
$ws.Range("T14").Value = 0.03046050470482393
$ws.Range("G15").Value = 0.27455
$ws.Range("H15").Value = 0.5491
$ws.Range("I15").Value = 0.3111729859610415
$ws.Range("J15").Value = 0.2314563074056745
$ws.Range("M15").Value = 13.324299
$ws.Range("N15").Value = 26.648598
$ws.Range("O15").Value = 0.2549418168249328
$ws.Range("P15").Value = 0.1857460205306503
$ws.Range("Q15").Value = 3.65818629045
$ws.Range("R15").Value = 14.6327451618
$ws.Range("S15").Value = 0.07933100638774725
$ws.Range("T15").Value = 0.04299208802732291
$ws.Range("G16").Value = 0.27455
$ws.Range("H16").Value = 0.5491
$ws.Range("I16").Value = 0.3111729859610415
$ws.Range("J16").Value = 0.2314563074056745
$ws.Range("M16").Value = 5.800427666666667
$ws.Range("N16").Value = 17.401283
$ws.Range("O16").Value = 0.110983066929195
$ws.Range("P16").Value = 0.1212903984433873
$ws.Range("Q16").Value = 1.592507415883333
$ws.Range("R16").Value = 9.555044495300001
$ws.Range("S16").Value = 0.03453493232747172
$ws.Range("T16").Value = 0.0280734277474694
$ws.Range("G17").Value = 0.2318183333333333
$ws.Range("H17").Value = 0.695455
$ws.Range("I17").Value = 0.2627412237619573
$ws.Range("J17").Value = 0.2931477804895526
$ws.Range("M17").Value = 7.004922666666666
$ws.Range("N17").Value = 21.014768
$ws.Range("O17").Value = 0.134029393318039
$ws.Range("P17").Value = 0.1464771065395205
$ws.Range("Q17").Value = 1.623869497715556
$ws.Range("R17").Value = 14.61482547944
$ws.Range("S17").Value = 0.03521504682045426
$ws.Range("T17").Value = 0.04293943867459216
$ws.Range("G18").Value = 0.2318183333333333
$ws.Range("H18").Value = 0.695455
$ws.Range("I18").Value = 0.2627412237619573
$ws.Range("J18").Value = 0.2931477804895526
$ws.Range("O18").Value = 0.3796257919253833
$ws.Range("P18").Value = 0.4148827819958515
$ws.Range("Q18").Value = 4.599459333453334
$ws.Range("R18").Value = 41.39513400108
$ws.Range("S18").Value = 0.09974334514207736
$ws.Range("T18").Value = 0.1216219667054148
$ws.Range("G19").Value = 0.2318183333333333
$ws.Range("H19").Value = 0.695455
$ws.Range("I19").Value = 0.2627412237619573
$ws.Range("J19").Value = 0.2931477804895526
$ws.Range("M19").Value = 6.293636666666667
$ws.Range("N19").Value = 18.88091
$ws.Range("O19").Value = 0.12041993100245
$ws.Range("P19").Value = 0.1316036924905903
$ws.Range("Q19").Value = 1.458980362672222
$ws.Range("R19").Value = 13.13082326405
$ws.Range("S19").Value = 0.03163928003691419
$ws.Range("T19").Value = 0.03857933035784616
$ws.Range("G20").Value = 0.2318183333333333
$ws.Range("H20").Value = 0.695455
$ws.Range("I20").Value = 0.2627412237619573
$ws.Range("J20").Value = 0.2931477804895526
$ws.Range("M20").Value = 13.324299
$ws.Range("N20").Value = 26.648598
$ws.Range("O20").Value = 0.2549418168249328
$ws.Range("P20").Value = 0.1857460205306503
$ws.Range("Q20").Value = 3.088816787015
$ws.Range("R20").Value = 18.53290072209
$ws.Range("S20").Value = 0.06698372494067961
$ws.Range("T20").Value = 0.054451033653327
$ws.Range("G21").Value = 0.2318183333333333
$ws.Range("H21").Value = 0.695455
$ws.Range("I21").Value = 0.2627412237619573
$ws.Range("J21").Value = 0.2931477804895526
$ws.Range("M21").Value = 5.800427666666667
$ws.Range("N21").Value = 17.401283
$ws.Range("O21").Value = 0.110983066929195
$ws.Range("P21").Value = 0.1212903984433873
$ws.Range("Q21").Value = 1.344645474307222
$ws.Range("R21").Value = 12.101809268765
$ws.Range("S21").Value = 0.0291598268218319
$ws.Range("T21").Value = 0.03555601109837249
